# RPA datasets push 2024-06-17
# Insert a new IPO-pipeline record ("라메디텍", underwriter "대신") as the new
# second row of the "01_IB전략컨설팅부" sheet, pushing all existing rows down
# by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push existing data rows (old row 2..13) down to row 3..14, and make row 2
# available for the new record.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:Y2").ClearFormats()

# Columns A-E, N, O and Y hold text (dates are stored as plain text strings,
# not real dates, and the odds/percentage columns are text too) - force text
# number-formatting before assigning so Excel doesn't auto-convert them to
# date serials / percentages.
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("N2:O2").NumberFormat = "@"
$ws.Range("Y2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-05-27"
$ws.Range("B2").Value = "2024-05-31"
$ws.Range("C2").Value = "2024-06-17"
$ws.Range("D2").Value = "대신"
$ws.Range("E2").Value = "라메디텍"
$ws.Range("F2").Value = 1298000
$ws.Range("G2").Value = 1298000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 10400
$ws.Range("J2").Value = 12700
$ws.Range("K2").Value = 8650735
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 16000
$ws.Range("N2").Value = "1115.44:1"
$ws.Range("O2").Value = "9.93%"
$ws.Range("P2").Value = 2055746777
$ws.Range("Q2").Value = 2918221978
$ws.Range("R2").Value = 979078233
$ws.Range("S2").Value = -3343774083
$ws.Range("T2").Value = -3525649863
$ws.Range("U2").Value = -1713494359
$ws.Range("V2").Value = -4430074915
$ws.Range("W2").Value = -8304699942
$ws.Range("X2").Value = -1627684107
$ws.Range("Y2").Value = "초소형 레이저 의료기기 및 미용기기"

# Drop the text number-format left behind by the assignments above so the
# new row matches the unstyled look of the other data rows.
$ws.Range("A2:Y2").ClearFormats()
